$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.195759773254395
$ws.Range("B1").Value = 2.177659273147583
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.142472743988037
$ws.Range("E1").Value = 1.212908148765564
